$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGD7QOAF4%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T210005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQCRejbTTyWT4e3boBavHXYIvkjvCtaFS3jBdbbm7Gsl7gIgeG4UAbVLoJViHY%2BuTduSdUigD0cPRYm7snA5tl1skKsq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJwrYtmSAiX2xn4yMSq6A5%2BmL2lVibYHPxEwWlivv1bGMMfRa7od1KkxJOIg4%2Fb%2Foab9qc9%2Bk%2BMU3333mFIPbp3FSn7vmmkt4tV%2FrzmUdvHxgFQ9atyCzRB3m4AyDYRLLBdTlM6Kd1CYU2jnQYlI3YIveo8WNVi90tysLuLRdJjr2HL8Dtv3V0csmJ44jMj3rp%2BfHxeMJ%2BqWsDB4S3xejUJkxhczAVMTXhb8MddG8XMdQesA6ZLWHwpXAVWaiL160d0Y0YRdRUYpX2z1XiYgRxSY8v4SrhgsQabMbaZt34Ud21wJn7q%2Fp%2B733u%2FQfgjuNMgmR9Qgrsz3bBuPgWq1dMaQB4zj51gmenftGZEoi8CZxpwzgMDodVAUQq%2BMPKLgqhj3HVZNgqS6ALBnNjLu1KR9GtfW0LC2knAMRE5Bk65AZIm5Nt8QQi6w8XBW0N3z3HYZETdoIf%2Fh0VE%2FhJ0pTBZJ4Llm6gAfwXU28eZoK00tnAATrD6QGyZmqrUO0qJf%2F%2BlqLATHgFBDlY4kE6Qv36yyzPW91S74TmefNaQ5D849aXlJoRXCnITMtWYP51pW2GaVW5Fc1xhgDX9D067FENWcMdONtJTuHrcwxMmlvwY6pQEihnYxKg%2FkxjhBDOTS9SfdMd0lXkv4Q7OSbJGYVsEaIThBXDIJ3y9ICbnQq4LomPYgtCNjdYFHoHs9aVqGe0UgiXo3unKAl7w7N1Qvm8P6KKSsTKK%2Blad%2BmIc1%2BhsZm%2Fq92c03wcJRGuxqHic7F7vKbP32spa5H2VHhPSZ3PwMVTJ%2BKEclVgJnp5oOyRjSQyJ%2FyR5Dw%2FS4YXItn5zOcJh6xDlQeCw%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2005%3A00%3A05%20GMT&X-Amz-Signature=a69da15d4dedfec174bd4bb825eb55df63792b53a309ad57fab1246d84c6c84a'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGD7QOAF4%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T210005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQCRejbTTyWT4e3boBavHXYIvkjvCtaFS3jBdbbm7Gsl7gIgeG4UAbVLoJViHY%2BuTduSdUigD0cPRYm7snA5tl1skKsq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJwrYtmSAiX2xn4yMSq6A5%2BmL2lVibYHPxEwWlivv1bGMMfRa7od1KkxJOIg4%2Fb%2Foab9qc9%2Bk%2BMU3333mFIPbp3FSn7vmmkt4tV%2FrzmUdvHxgFQ9atyCzRB3m4AyDYRLLBdTlM6Kd1CYU2jnQYlI3YIveo8WNVi90tysLuLRdJjr2HL8Dtv3V0csmJ44jMj3rp%2BfHxeMJ%2BqWsDB4S3xejUJkxhczAVMTXhb8MddG8XMdQesA6ZLWHwpXAVWaiL160d0Y0YRdRUYpX2z1XiYgRxSY8v4SrhgsQabMbaZt34Ud21wJn7q%2Fp%2B733u%2FQfgjuNMgmR9Qgrsz3bBuPgWq1dMaQB4zj51gmenftGZEoi8CZxpwzgMDodVAUQq%2BMPKLgqhj3HVZNgqS6ALBnNjLu1KR9GtfW0LC2knAMRE5Bk65AZIm5Nt8QQi6w8XBW0N3z3HYZETdoIf%2Fh0VE%2FhJ0pTBZJ4Llm6gAfwXU28eZoK00tnAATrD6QGyZmqrUO0qJf%2F%2BlqLATHgFBDlY4kE6Qv36yyzPW91S74TmefNaQ5D849aXlJoRXCnITMtWYP51pW2GaVW5Fc1xhgDX9D067FENWcMdONtJTuHrcwxMmlvwY6pQEihnYxKg%2FkxjhBDOTS9SfdMd0lXkv4Q7OSbJGYVsEaIThBXDIJ3y9ICbnQq4LomPYgtCNjdYFHoHs9aVqGe0UgiXo3unKAl7w7N1Qvm8P6KKSsTKK%2Blad%2BmIc1%2BhsZm%2Fq92c03wcJRGuxqHic7F7vKbP32spa5H2VHhPSZ3PwMVTJ%2BKEclVgJnp5oOyRjSQyJ%2FyR5Dw%2FS4YXItn5zOcJh6xDlQeCw%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2005%3A00%3A05%20GMT&X-Amz-Signature=b9e1be06d7509036caf3803c551d09524721ae1abd3705f127ede9e5086d6246'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGD7QOAF4%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T210005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQCRejbTTyWT4e3boBavHXYIvkjvCtaFS3jBdbbm7Gsl7gIgeG4UAbVLoJViHY%2BuTduSdUigD0cPRYm7snA5tl1skKsq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJwrYtmSAiX2xn4yMSq6A5%2BmL2lVibYHPxEwWlivv1bGMMfRa7od1KkxJOIg4%2Fb%2Foab9qc9%2Bk%2BMU3333mFIPbp3FSn7vmmkt4tV%2FrzmUdvHxgFQ9atyCzRB3m4AyDYRLLBdTlM6Kd1CYU2jnQYlI3YIveo8WNVi90tysLuLRdJjr2HL8Dtv3V0csmJ44jMj3rp%2BfHxeMJ%2BqWsDB4S3xejUJkxhczAVMTXhb8MddG8XMdQesA6ZLWHwpXAVWaiL160d0Y0YRdRUYpX2z1XiYgRxSY8v4SrhgsQabMbaZt34Ud21wJn7q%2Fp%2B733u%2FQfgjuNMgmR9Qgrsz3bBuPgWq1dMaQB4zj51gmenftGZEoi8CZxpwzgMDodVAUQq%2BMPKLgqhj3HVZNgqS6ALBnNjLu1KR9GtfW0LC2knAMRE5Bk65AZIm5Nt8QQi6w8XBW0N3z3HYZETdoIf%2Fh0VE%2FhJ0pTBZJ4Llm6gAfwXU28eZoK00tnAATrD6QGyZmqrUO0qJf%2F%2BlqLATHgFBDlY4kE6Qv36yyzPW91S74TmefNaQ5D849aXlJoRXCnITMtWYP51pW2GaVW5Fc1xhgDX9D067FENWcMdONtJTuHrcwxMmlvwY6pQEihnYxKg%2FkxjhBDOTS9SfdMd0lXkv4Q7OSbJGYVsEaIThBXDIJ3y9ICbnQq4LomPYgtCNjdYFHoHs9aVqGe0UgiXo3unKAl7w7N1Qvm8P6KKSsTKK%2Blad%2BmIc1%2BhsZm%2Fq92c03wcJRGuxqHic7F7vKbP32spa5H2VHhPSZ3PwMVTJ%2BKEclVgJnp5oOyRjSQyJ%2FyR5Dw%2FS4YXItn5zOcJh6xDlQeCw%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2005%3A00%3A05%20GMT&X-Amz-Signature=3c859ee5b8eee385d86699d545b88b4e4585f35ee7b10e22e31cca96e61b0033'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGD7QOAF4%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T210005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQCRejbTTyWT4e3boBavHXYIvkjvCtaFS3jBdbbm7Gsl7gIgeG4UAbVLoJViHY%2BuTduSdUigD0cPRYm7snA5tl1skKsq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJwrYtmSAiX2xn4yMSq6A5%2BmL2lVibYHPxEwWlivv1bGMMfRa7od1KkxJOIg4%2Fb%2Foab9qc9%2Bk%2BMU3333mFIPbp3FSn7vmmkt4tV%2FrzmUdvHxgFQ9atyCzRB3m4AyDYRLLBdTlM6Kd1CYU2jnQYlI3YIveo8WNVi90tysLuLRdJjr2HL8Dtv3V0csmJ44jMj3rp%2BfHxeMJ%2BqWsDB4S3xejUJkxhczAVMTXhb8MddG8XMdQesA6ZLWHwpXAVWaiL160d0Y0YRdRUYpX2z1XiYgRxSY8v4SrhgsQabMbaZt34Ud21wJn7q%2Fp%2B733u%2FQfgjuNMgmR9Qgrsz3bBuPgWq1dMaQB4zj51gmenftGZEoi8CZxpwzgMDodVAUQq%2BMPKLgqhj3HVZNgqS6ALBnNjLu1KR9GtfW0LC2knAMRE5Bk65AZIm5Nt8QQi6w8XBW0N3z3HYZETdoIf%2Fh0VE%2FhJ0pTBZJ4Llm6gAfwXU28eZoK00tnAATrD6QGyZmqrUO0qJf%2F%2BlqLATHgFBDlY4kE6Qv36yyzPW91S74TmefNaQ5D849aXlJoRXCnITMtWYP51pW2GaVW5Fc1xhgDX9D067FENWcMdONtJTuHrcwxMmlvwY6pQEihnYxKg%2FkxjhBDOTS9SfdMd0lXkv4Q7OSbJGYVsEaIThBXDIJ3y9ICbnQq4LomPYgtCNjdYFHoHs9aVqGe0UgiXo3unKAl7w7N1Qvm8P6KKSsTKK%2Blad%2BmIc1%2BhsZm%2Fq92c03wcJRGuxqHic7F7vKbP32spa5H2VHhPSZ3PwMVTJ%2BKEclVgJnp5oOyRjSQyJ%2FyR5Dw%2FS4YXItn5zOcJh6xDlQeCw%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2005%3A00%3A05%20GMT&X-Amz-Signature=ed8de039cda36d001557be5b84f7311a3fdefec290c8554a29b891ff9b41f1c0'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGD7QOAF4%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T210005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQCRejbTTyWT4e3boBavHXYIvkjvCtaFS3jBdbbm7Gsl7gIgeG4UAbVLoJViHY%2BuTduSdUigD0cPRYm7snA5tl1skKsq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJwrYtmSAiX2xn4yMSq6A5%2BmL2lVibYHPxEwWlivv1bGMMfRa7od1KkxJOIg4%2Fb%2Foab9qc9%2Bk%2BMU3333mFIPbp3FSn7vmmkt4tV%2FrzmUdvHxgFQ9atyCzRB3m4AyDYRLLBdTlM6Kd1CYU2jnQYlI3YIveo8WNVi90tysLuLRdJjr2HL8Dtv3V0csmJ44jMj3rp%2BfHxeMJ%2BqWsDB4S3xejUJkxhczAVMTXhb8MddG8XMdQesA6ZLWHwpXAVWaiL160d0Y0YRdRUYpX2z1XiYgRxSY8v4SrhgsQabMbaZt34Ud21wJn7q%2Fp%2B733u%2FQfgjuNMgmR9Qgrsz3bBuPgWq1dMaQB4zj51gmenftGZEoi8CZxpwzgMDodVAUQq%2BMPKLgqhj3HVZNgqS6ALBnNjLu1KR9GtfW0LC2knAMRE5Bk65AZIm5Nt8QQi6w8XBW0N3z3HYZETdoIf%2Fh0VE%2FhJ0pTBZJ4Llm6gAfwXU28eZoK00tnAATrD6QGyZmqrUO0qJf%2F%2BlqLATHgFBDlY4kE6Qv36yyzPW91S74TmefNaQ5D849aXlJoRXCnITMtWYP51pW2GaVW5Fc1xhgDX9D067FENWcMdONtJTuHrcwxMmlvwY6pQEihnYxKg%2FkxjhBDOTS9SfdMd0lXkv4Q7OSbJGYVsEaIThBXDIJ3y9ICbnQq4LomPYgtCNjdYFHoHs9aVqGe0UgiXo3unKAl7w7N1Qvm8P6KKSsTKK%2Blad%2BmIc1%2BhsZm%2Fq92c03wcJRGuxqHic7F7vKbP32spa5H2VHhPSZ3PwMVTJ%2BKEclVgJnp5oOyRjSQyJ%2FyR5Dw%2FS4YXItn5zOcJh6xDlQeCw%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2005%3A00%3A05%20GMT&X-Amz-Signature=48db557d71cb34a0449d2a876f38a47af8b2c405ad8ce789a7112ad2fde15d0d'
